$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'44.756.58"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +3.83%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.420.50"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.95%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.01%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'317.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +4.63%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'101.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +6.30%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.513"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +1.98%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.09%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.528"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +10.03%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'35.34"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +2.88%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0800"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +1.57%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -1.66%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'18.59"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.95%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'6.89"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +2.15%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'2.798.24"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +2.16%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'2.457.84"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +3.82%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D18").Value = "'44.577.92"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +3.39%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'12.24"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +1.94%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  +1.15%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.0₃0919"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +3.49%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'68.57"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.83%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'242.35"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +2.82%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +3.48%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +2.24%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -0.02%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'25.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +2.91%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -3.59%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'9.50"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +1.86%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'33.33"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +3.27%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'48.35"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.88%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.126"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +14.68%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'19.58"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +11.07%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +2.83%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +0.21%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.0765"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +4.81%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +2.34%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'4.46"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +3.33%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'126.62"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.38%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -0.34%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.110"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +1.75%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'2.18"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -3.83%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'20.94"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.72%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.0289"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +3.58%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'1.938.84"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.43%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -0.68%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'2.93"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +7.02%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'9.15"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.69%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  +15.59%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'75.42"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +5.50%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'53.64"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +4.89%  "
$ws.Range("E51").Style = "Normal"
